$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# Remove "area" column (old G1), shift lat/lng left, add "lnk" and "bw" headers
$ws.Range("G1").Value = "lat"
$ws.Range("H1").Value = "lng"
$ws.Range("C2").Value = "custom street 1st"
$ws.Range("I1").Value = "lnk"
$ws.Range("I2").Value = "Astinet"
$ws.Range("J1").Value = "bw"
$ws.Range("J2").Value = "2 MB"

# --- Data row (row 2) ---
$ws.Range("A2").Value = "LOC123"
$ws.Range("B2").Value = "location 123"

# --- Column widths (match "best fit" widths computed by Excel for the new data) ---
$ws.Columns.Item(2).ColumnWidth = 8.666666666666666
$ws.Columns.Item(3).ColumnWidth = 15.5

# --- Selection ---
$ws.Range("G1").Select()
